$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New "Save" column header in H1, using same formatting as the other header cells (e.g. G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# New "Save" values for the data rows
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
